# Pending Questions.docx -- add a new pending question to the list.
#
# 1. The existing hidden "_GoBack" bookmark (left over from the previous
#    edit, sitting inside the "Would you rather us ... implement a GUI or
#    console application?" paragraph) is removed. Word keeps only one
#    "_GoBack" bookmark at a time, relocating it to the site of the most
#    recent edit.
# 2. Two new paragraphs are appended after the "What type of report..."
#    question: a blank ListParagraph-styled spacer, then a new numbered
#    list item holding the new question. The "_GoBack" bookmark is
#    re-created at the end of that new question's text, marking it as the
#    location of the latest edit -- exactly as Word itself would do.

$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$lastQuestion = $d.Paragraphs.Item(11)

# Grow the document by two empty paragraphs right after the last question,
# cloning its list-numbering formatting; we then overwrite each paragraph's
# contents/formatting in place below.
$lastQuestion.Range.InsertParagraphAfter()
$spacerParagraph = $d.Paragraphs.Item(12)
$spacerParagraph.Range.InsertParagraphAfter()

# Paragraph 1 of 2: a blank "ListParagraph" spacer line (no numbering, no
# explicit spacing override) -- matches the blank line that separates every
# other question in the list.
$spacerParagraph = $d.Paragraphs.Item(12)
$spacerXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$spacerParagraph.Range.InsertXML($spacerXml)

# Paragraph 2 of 2: the new numbered question itself, with the "_GoBack"
# bookmark wrapping the end of its text (an empty, collapsed bookmark span)
# just like Word leaves behind after typing new text.
$newQuestionParagraph = $d.Paragraphs.Item(13)
$newQuestionXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Does the procedure that creates weekly reports for members need to able to run automatically every week? Or is it triggered manually?</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newQuestionParagraph.Range.InsertXML($newQuestionXml)
